$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PowerPlantsPerformance")
$ws.Activate()
$ws.Range("H29").Select()
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 1
